# Add settings, TDDied the organizations list, use FSharp data for parsing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- AP7: was a plain text cell, now pulls the "ksiegowanie" phone label
#     straight off the header (P1) via a formula instead of a duplicated
#     literal string.
$ws.Range("AP7").Formula = "=+P1"

# --- Selection / view: user had scrolled over to AP7; new workbook opens
#     back at the top with L4 selected.
$ws.Range("L4").Select()

# --- Unhide the S:U helper columns (regon / krs / powiat helper block)
$ws.Range("S1:U1").EntireColumn.Hidden = $false

# --- Row 13 used to be a single SUM() footer cell; it is replaced with a
#     row of section headings that TDD/describe the column groups below,
#     and the beneficiary total becomes a label instead of the SUM formula.
$ws.Range("AC13").Value = "Beneficjenci"
$ws.Range("B13").Value = "Numery"
$ws.Range("H13").Value = "Dane adresowe"
$ws.Range("N13").Value = "Dane adresowe ksiegowosci"
$ws.Range("Q13").Value = "Kontakty"
$ws.Range("AE13").Value = "Źródła żywności"
$ws.Range("AI13").Value = "Warunki udzielania pomocy żywnościowej"
$ws.Range("AQ13").Value = "Dokumentacja"
